$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Summary" block (rows 45-49) ---------------------------------
# Header row
# (String insertion order below is chosen to reproduce the shared-string
#  table ordering that the original author's edit produced.)

$ws.Range("A47").Value = "Functional in IMGT and Digger"
$ws.Range("B47").Formula = "=B19"
$ws.Range("C47").Formula = "=B32"
$ws.Range("D47").Formula = "=B47+C47"

$ws.Range("A48").Value = "Functional in IMGT only"
$ws.Range("B48").Formula = "=C19"
$ws.Range("C48").Formula = "=C32"

$ws.Range("A49").Value = "Functional in Digger only"
$ws.Range("B49").Formula = "=D19"
$ws.Range("C49").Formula = "=D32"

# Shared formula across D48:D49
$ws.Range("D48:D49").Formula = "=B48+C48"

$ws.Range("A46").Value = "Summary"
$ws.Range("A46").Font.Italic = $true

# --- New "Analysis of differences" block (rows 50-55) ------------------
$ws.Range("A50").Value = "Analysis of differences"
$ws.Range("A50").Font.Italic = $true

$ws.Range("A51").Value = "non-annotation differences"
$ws.Range("B51").Formula = "=E19"
$ws.Range("C51").Formula = "=E32"

$ws.Range("A53").Value = "LEADER"
$ws.Range("B53").Formula = "=G19"
$ws.Range("C53").Formula = "=G32"

$ws.Range("A52").Value = "STOP-CODON"
$ws.Range("B52").Formula = "=F19"
$ws.Range("C52").Formula = "=F32"

$ws.Range("A54").Value = "RSS"
$ws.Range("B54").Formula = "=H19"
$ws.Range("C54").Formula = "=H32"

$ws.Range("A55").Value = "Sequence not identified"
$ws.Range("B55").Formula = "=I19"
$ws.Range("C55").Formula = "=I32"

# Shared formula across D51:D55
$ws.Range("D51:D55").Formula = "=B51+C51"

# Column headers above the Summary block (row 45)
$ws.Range("B45").Value = "Human IG, TR"
$ws.Range("C45").Value = "Rhesus Macaque IG"
$ws.Range("D45").Value = "TOTAL"

# --- Misc fix noted in the diff: E32 formula rewritten explicitly ------
$ws.Range("E32").Formula = "=SUM(E23:E30)"

# --- Sheet cosmetics -----------------------------------------------------
$ws.Columns("A").ColumnWidth = 47.6

$ws.Range("B36").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
